$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
# Clear old values
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
# Set new values
$ws.Range("G2").Value = 8
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 2

# --- Row 3 ---
$ws.Range("U3").ClearContents()
$ws.Range("Y3").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AA3").ClearContents()
$ws.Range("V3").Value = 3
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 8
$ws.Range("AB3").Value = 8
$ws.Range("AC3").Value = 8
$ws.Range("AD3").Value = 8
$ws.Range("AE3").Value = 8

# --- Row 4 ---
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("R4").ClearContents()
$ws.Range("S4").ClearContents()
$ws.Range("T4").ClearContents()
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 8
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 8
$ws.Range("P4").Value = 8
$ws.Range("Q4").Value = 8
$ws.Range("U4").Value = 8
$ws.Range("V4").Value = 5
$ws.Range("W4").Value = 1

# --- Sheet view changes: zoom and selection ---
$ws.Select() | Out-Null
$excel.ActiveWindow.Zoom = 130
$ws.Range("AH6").Select() | Out-Null
